# This table is a "skimr"-style data-summary export. The edit adds a new
# "character" row describing a variable named "mujer" right after the
# existing "JHOGAR" row, which pushes every subsequent "numeric" variable
# row down by one, and slightly updates their recomputed summary
# statistics (mean/sd/percentiles/histogram) to reflect the updated
# dataset.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 3 for the new "mujer" character variable
# (skim_type=character). This shifts the previous rows 3-17 (the numeric
# variable summaries) down to rows 4-18.
$ws.Rows.Item(3).Insert()

# Row 2: JHOGAR
$ws.Cells.Item(2,1).Value2 = "character"
$ws.Cells.Item(2,2).Value2 = "JHOGAR"
$ws.Cells.Item(2,3).Value2 = 0
$ws.Cells.Item(2,4).Value2 = 1
$ws.Cells.Item(2,5).Value2 = 1
$ws.Cells.Item(2,6).Value2 = 1
$ws.Cells.Item(2,7).Value2 = 0
$ws.Cells.Item(2,8).Value2 = 2
$ws.Cells.Item(2,9).Value2 = 0

# Row 3: mujer
$ws.Cells.Item(3,1).Value2 = "character"
$ws.Cells.Item(3,2).Value2 = "mujer"
$ws.Cells.Item(3,3).Value2 = 0
$ws.Cells.Item(3,4).Value2 = 1
$ws.Cells.Item(3,5).Value2 = 1
$ws.Cells.Item(3,6).Value2 = 1
$ws.Cells.Item(3,7).Value2 = 0
$ws.Cells.Item(3,8).Value2 = 2
$ws.Cells.Item(3,9).Value2 = 0

# Row 4: age
$ws.Cells.Item(4,1).Value2 = "numeric"
$ws.Cells.Item(4,2).Value2 = "age"
$ws.Cells.Item(4,3).Value2 = 0
$ws.Cells.Item(4,4).Value2 = 1
$ws.Cells.Item(4,10).Value2 = 39.0766129032258
$ws.Cells.Item(4,11).Value2 = 13.1065806455642
$ws.Cells.Item(4,12).Value2 = 19
$ws.Cells.Item(4,13).Value2 = 28
$ws.Cells.Item(4,14).Value2 = 37
$ws.Cells.Item(4,15).Value2 = 49
$ws.Cells.Item(4,16).Value2 = 91
$ws.Cells.Item(4,17).Value2 = "▇▆▅▁▁"

# Row 5: cuentaPropia
$ws.Cells.Item(5,1).Value2 = "numeric"
$ws.Cells.Item(5,2).Value2 = "cuentaPropia"
$ws.Cells.Item(5,3).Value2 = 0
$ws.Cells.Item(5,4).Value2 = 1
$ws.Cells.Item(5,10).Value2 = 0.298455440131219
$ws.Cells.Item(5,11).Value2 = 0.457596002040128
$ws.Cells.Item(5,12).Value2 = 0
$ws.Cells.Item(5,13).Value2 = 0
$ws.Cells.Item(5,14).Value2 = 0
$ws.Cells.Item(5,15).Value2 = 1
$ws.Cells.Item(5,16).Value2 = 1
$ws.Cells.Item(5,17).Value2 = "▇▁▁▁▃"

# Row 6: formal
$ws.Cells.Item(6,1).Value2 = "numeric"
$ws.Cells.Item(6,2).Value2 = "formal"
$ws.Cells.Item(6,3).Value2 = 0
$ws.Cells.Item(6,4).Value2 = 1
$ws.Cells.Item(6,10).Value2 = 0.605795516675779
$ws.Cells.Item(6,11).Value2 = 0.488695846779939
$ws.Cells.Item(6,12).Value2 = 0
$ws.Cells.Item(6,13).Value2 = 0
$ws.Cells.Item(6,14).Value2 = 1
$ws.Cells.Item(6,15).Value2 = 1
$ws.Cells.Item(6,16).Value2 = 1
$ws.Cells.Item(6,17).Value2 = "▅▁▁▁▇"

# Row 7: hoursWorkUsual
$ws.Cells.Item(7,1).Value2 = "numeric"
$ws.Cells.Item(7,2).Value2 = "hoursWorkUsual"
$ws.Cells.Item(7,3).Value2 = 0
$ws.Cells.Item(7,4).Value2 = 1
$ws.Cells.Item(7,10).Value2 = 47.2450792782941
$ws.Cells.Item(7,11).Value2 = 15.0131265229188
$ws.Cells.Item(7,12).Value2 = 1
$ws.Cells.Item(7,13).Value2 = 40
$ws.Cells.Item(7,14).Value2 = 48
$ws.Cells.Item(7,15).Value2 = 50
$ws.Cells.Item(7,16).Value2 = 130
$ws.Cells.Item(7,17).Value2 = "▁▇▂▁▁"

# Row 8: maxEducLevel
$ws.Cells.Item(8,1).Value2 = "numeric"
$ws.Cells.Item(8,2).Value2 = "maxEducLevel"
$ws.Cells.Item(8,3).Value2 = 0
$ws.Cells.Item(8,4).Value2 = 1
$ws.Cells.Item(8,10).Value2 = 5.94594040459267
$ws.Cells.Item(8,11).Value2 = 1.21273624318973
$ws.Cells.Item(8,12).Value2 = 0
$ws.Cells.Item(8,13).Value2 = 5
$ws.Cells.Item(8,14).Value2 = 6
$ws.Cells.Item(8,15).Value2 = 7
$ws.Cells.Item(8,16).Value2 = 7
$ws.Cells.Item(8,17).Value2 = "▁▁▂▁▇"

# Row 9: ocu
$ws.Cells.Item(9,1).Value2 = "numeric"
$ws.Cells.Item(9,2).Value2 = "ocu"
$ws.Cells.Item(9,3).Value2 = 0
$ws.Cells.Item(9,4).Value2 = 1
$ws.Cells.Item(9,10).Value2 = 1
$ws.Cells.Item(9,11).Value2 = 0
$ws.Cells.Item(9,12).Value2 = 1
$ws.Cells.Item(9,13).Value2 = 1
$ws.Cells.Item(9,14).Value2 = 1
$ws.Cells.Item(9,15).Value2 = 1
$ws.Cells.Item(9,16).Value2 = 1
$ws.Cells.Item(9,17).Value2 = "▁▁▇▁▁"

# Row 10: oficio
$ws.Cells.Item(10,1).Value2 = "numeric"
$ws.Cells.Item(10,2).Value2 = "oficio"
$ws.Cells.Item(10,3).Value2 = 0
$ws.Cells.Item(10,4).Value2 = 1
$ws.Cells.Item(10,10).Value2 = 50.1713367960634
$ws.Cells.Item(10,11).Value2 = 28.0737206019987
$ws.Cells.Item(10,12).Value2 = 1
$ws.Cells.Item(10,13).Value2 = 33
$ws.Cells.Item(10,14).Value2 = 45
$ws.Cells.Item(10,15).Value2 = 74
$ws.Cells.Item(10,16).Value2 = 99
$ws.Cells.Item(10,17).Value2 = "▃▃▇▁▅"

# Row 11: estrato1
$ws.Cells.Item(11,1).Value2 = "numeric"
$ws.Cells.Item(11,2).Value2 = "estrato1"
$ws.Cells.Item(11,3).Value2 = 0
$ws.Cells.Item(11,4).Value2 = 1
$ws.Cells.Item(11,10).Value2 = 2.5247402952433
$ws.Cells.Item(11,11).Value2 = 0.988936268493968
$ws.Cells.Item(11,12).Value2 = 1
$ws.Cells.Item(11,13).Value2 = 2
$ws.Cells.Item(11,14).Value2 = 2
$ws.Cells.Item(11,15).Value2 = 3
$ws.Cells.Item(11,16).Value2 = 6
$ws.Cells.Item(11,17).Value2 = "▇▆▁▁▁"

# Row 12: informal
$ws.Cells.Item(12,1).Value2 = "numeric"
$ws.Cells.Item(12,2).Value2 = "informal"
$ws.Cells.Item(12,3).Value2 = 0
$ws.Cells.Item(12,4).Value2 = 1
$ws.Cells.Item(12,10).Value2 = 0.394204483324221
$ws.Cells.Item(12,11).Value2 = 0.488695846779939
$ws.Cells.Item(12,12).Value2 = 0
$ws.Cells.Item(12,13).Value2 = 0
$ws.Cells.Item(12,14).Value2 = 0
$ws.Cells.Item(12,15).Value2 = 1
$ws.Cells.Item(12,16).Value2 = 1
$ws.Cells.Item(12,17).Value2 = "▇▁▁▁▅"

# Row 13: p6050
$ws.Cells.Item(13,1).Value2 = "numeric"
$ws.Cells.Item(13,2).Value2 = "p6050"
$ws.Cells.Item(13,3).Value2 = 0
$ws.Cells.Item(13,4).Value2 = 1
$ws.Cells.Item(13,10).Value2 = 2.18944778567523
$ws.Cells.Item(13,11).Value2 = 1.79376975108388
$ws.Cells.Item(13,12).Value2 = 1
$ws.Cells.Item(13,13).Value2 = 1
$ws.Cells.Item(13,14).Value2 = 2
$ws.Cells.Item(13,15).Value2 = 3
$ws.Cells.Item(13,16).Value2 = 9
$ws.Cells.Item(13,17).Value2 = "▇▂▁▁▁"

# Row 14: relab
$ws.Cells.Item(14,1).Value2 = "numeric"
$ws.Cells.Item(14,2).Value2 = "relab"
$ws.Cells.Item(14,3).Value2 = 0
$ws.Cells.Item(14,4).Value2 = 1
$ws.Cells.Item(14,10).Value2 = 2.14454620010935
$ws.Cells.Item(14,11).Value2 = 1.45685647304309
$ws.Cells.Item(14,12).Value2 = 1
$ws.Cells.Item(14,13).Value2 = 1
$ws.Cells.Item(14,14).Value2 = 1
$ws.Cells.Item(14,15).Value2 = 4
$ws.Cells.Item(14,16).Value2 = 9
$ws.Cells.Item(14,17).Value2 = "▇▅▁▁▁"

# Row 15: sex
$ws.Cells.Item(15,1).Value2 = "numeric"
$ws.Cells.Item(15,2).Value2 = "sex"
$ws.Cells.Item(15,3).Value2 = 0
$ws.Cells.Item(15,4).Value2 = 1
$ws.Cells.Item(15,10).Value2 = 0.526038819026791
$ws.Cells.Item(15,11).Value2 = 0.499338583097074
$ws.Cells.Item(15,12).Value2 = 0
$ws.Cells.Item(15,13).Value2 = 0
$ws.Cells.Item(15,14).Value2 = 1
$ws.Cells.Item(15,15).Value2 = 1
$ws.Cells.Item(15,16).Value2 = 1
$ws.Cells.Item(15,17).Value2 = "▇▁▁▁▇"

# Row 16: sizeFirm
$ws.Cells.Item(16,1).Value2 = "numeric"
$ws.Cells.Item(16,2).Value2 = "sizeFirm"
$ws.Cells.Item(16,3).Value2 = 0
$ws.Cells.Item(16,4).Value2 = 1
$ws.Cells.Item(16,10).Value2 = 3.21480317113177
$ws.Cells.Item(16,11).Value2 = 1.65477695271946
$ws.Cells.Item(16,12).Value2 = 1
$ws.Cells.Item(16,13).Value2 = 2
$ws.Cells.Item(16,14).Value2 = 4
$ws.Cells.Item(16,15).Value2 = 5
$ws.Cells.Item(16,16).Value2 = 5
$ws.Cells.Item(16,17).Value2 = "▅▃▂▃▇"

# Row 17: wap
$ws.Cells.Item(17,1).Value2 = "numeric"
$ws.Cells.Item(17,2).Value2 = "wap"
$ws.Cells.Item(17,3).Value2 = 0
$ws.Cells.Item(17,4).Value2 = 1
$ws.Cells.Item(17,10).Value2 = 1
$ws.Cells.Item(17,11).Value2 = 0
$ws.Cells.Item(17,12).Value2 = 1
$ws.Cells.Item(17,13).Value2 = 1
$ws.Cells.Item(17,14).Value2 = 1
$ws.Cells.Item(17,15).Value2 = 1
$ws.Cells.Item(17,16).Value2 = 1
$ws.Cells.Item(17,17).Value2 = "▁▁▇▁▁"

# Row 18: y_total_m
$ws.Cells.Item(18,1).Value2 = "numeric"
$ws.Cells.Item(18,2).Value2 = "y_total_m"
$ws.Cells.Item(18,3).Value2 = 0
$ws.Cells.Item(18,4).Value2 = 1
$ws.Cells.Item(18,10).Value2 = 1626339.58406713
$ws.Cells.Item(18,11).Value2 = 2440279.41410928
$ws.Cells.Item(18,12).Value2 = 84
$ws.Cells.Item(18,13).Value2 = 800000
$ws.Cells.Item(18,14).Value2 = 996556.5
$ws.Cells.Item(18,15).Value2 = 1562500
$ws.Cells.Item(18,16).Value2 = 70000000
$ws.Cells.Item(18,17).Value2 = "▇▁▁▁▁"
